$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46022
$ws.Range("B2").Value = 100.4
$ws.Range("C2").Value = 85.36
$ws.Range("D2").Value = 81.54000000000001
$ws.Range("E2").Value = 76.58
$ws.Range("F2").Value = 75.16
$ws.Range("G2").Value = 77.39
$ws.Range("H2").Value = 83.59999999999999
$ws.Range("I2").Value = 95.08
$ws.Range("J2").Value = 105.94
$ws.Range("K2").Value = 102.26
$ws.Range("L2").Value = 97.34
$ws.Range("M2").Value = 84.5
$ws.Range("N2").Value = 83.05
$ws.Range("O2").Value = 85.17
$ws.Range("P2").Value = 87.59999999999999
$ws.Range("Q2").Value = 94.42
$ws.Range("R2").Value = 104.21
$ws.Range("S2").Value = 112.02
$ws.Range("T2").Value = 129.12
$ws.Range("U2").Value = 135.43
$ws.Range("V2").Value = 131.12
$ws.Range("W2").Value = 113.66
$ws.Range("X2").Value = 105.83
$ws.Range("Y2").Value = 100.95
$ws.Range("Z2").Value = 97.81999999999999
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 120.2
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 132.28
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 122.39
$ws.Range("AG2").Value = "1h-15h"
